# DeudoresPrueba.xlsx update
#
# 1) Remove five client rows that no longer belong in the debtor list:
#       CLIENTE PAOLA, EL RUBY, FRESIA, LA 13, SANDRA 20 DE JULIO
# 2) Insert a new client row for MAFE (right after LOS PAISANOS).
# 3) Correct a handful of date/value entries that were mis-keyed.
#
# Row numbers refer to the ORIGINAL sheet (rows 2..40 hold the data,
# row 1 is the header). Deletions are performed from the bottom up so
# that earlier row numbers stay valid while we work.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) delete obsolete rows (bottom-to-top so row numbers don't shift) ---
$ws.Rows(37).EntireRow.Delete()   # SANDRA 20 DE JULIO
$ws.Rows(22).EntireRow.Delete()   # LA 13
$ws.Rows(21).EntireRow.Delete()   # FRESIA
$ws.Rows(18).EntireRow.Delete()   # EL RUBY
$ws.Rows(13).EntireRow.Delete()   # CLIENTE PAOLA

# After the deletions above, the remaining 34 clients occupy rows 2..35,
# with LOS PAISANOS on row 21 and the two MERKA FRUVER DEXI rows on 22-23.

# --- 2) insert the new MAFE row right after LOS PAISANOS (row 21) ---
$ws.Rows(22).Insert()
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "MAFE"
$ws.Range("C22").Value = 46017
$ws.Range("D22").Value = 190000
$ws.Range("E22").Value = $false

# --- 3) fix up values that changed on the remaining/shifted rows ---

# CAMPO VERDE ZIPAQUIRA (row 7): corrected amount owed
$ws.Range("D7").Value = 425100

# MERKA FRUVER DEXI pair (now rows 23 & 24, pushed down one by the MAFE
# insert): the two entries were swapped
$ws.Range("C23").Value = 45988
$ws.Range("D23").Value = 15400
$ws.Range("C24").Value = 45995
$ws.Range("D24").Value = 339000

# PLAZA JESSICA (now row 30): corrected date and amount
$ws.Range("C30").Value = 46014
$ws.Range("D30").Value = 1655400

# --- renumber the Consecutivo column (A) sequentially for the final 35 rows ---
for ($r = 2; $r -le 36; $r++) {
    $ws.Range("A$r").Value = $r - 1
}
